# New crime data collected — update the 115th Precinct weekly CompStat
# report: header volume/date strings, and the Week-to-Date / 28-Day /
# Year-to-Date / 2-Year crime-complaint figures in rows 15-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 29   Number  44" -> "Volume 29   Number  45"
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "45"

# ---------------------------------------------------------------------
# Header: "Report Covering the Week  10/31/2022  Through  11/6/2022"
#      -> "Report Covering the Week  11/7/2022  Through  11/13/2022"
# ---------------------------------------------------------------------
$ws.Range("C9").Characters(27, 10).Text = "11/7/2022"
$ws.Range("C9").Characters(47, 9).Text = "11/13/2022"

# ---------------------------------------------------------------------
# Row 15 (Murder)
# ---------------------------------------------------------------------
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("J15").Value = 28
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = 40
$ws.Range("N15").Value = 20.689655172413

# ---------------------------------------------------------------------
# Row 16 (Rape)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 14.285714285714
$ws.Range("F16").Value = 37
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = 68.181818181818
$ws.Range("I16").Value = 269
$ws.Range("J16").Value = 214
$ws.Range("K16").Value = 25.700934579439
$ws.Range("L16").Value = 44.623655913978
$ws.Range("M16").Value = -12.944983818770
$ws.Range("N16").Value = -76.649305555555

# ---------------------------------------------------------------------
# Row 17 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -44.444444444444
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 7.407407407407
$ws.Range("I17").Value = 366
$ws.Range("J17").Value = 338
$ws.Range("K17").Value = 8.284023668639
$ws.Range("L17").Value = 36.059479553903
$ws.Range("M17").Value = 26.643598615917
$ws.Range("N17").Value = -2.659574468085

# ---------------------------------------------------------------------
# Row 18 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -29.411764705882
$ws.Range("I18").Value = 131
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 9.166666666666
$ws.Range("L18").Value = -11.486486486486
$ws.Range("M18").Value = -49.615384615384
$ws.Range("N18").Value = -92.396982008125

# ---------------------------------------------------------------------
# Row 19 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("E19").Value = 8.333333333333
$ws.Range("F19").Value = 86
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = 38.709677419354
$ws.Range("I19").Value = 854
$ws.Range("J19").Value = 478
$ws.Range("K19").Value = 78.661087866108
$ws.Range("L19").Value = 82.869379014989
$ws.Range("M19").Value = 98.143851508120
$ws.Range("N19").Value = -32.702915681639

# ---------------------------------------------------------------------
# Row 20 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 43.75
$ws.Range("I20").Value = 273
$ws.Range("J20").Value = 183
$ws.Range("K20").Value = 49.180327868852
$ws.Range("L20").Value = 51.666666666666
$ws.Range("M20").Value = 37.185929648241
$ws.Range("N20").Value = -86.128048780487

# ---------------------------------------------------------------------
# Row 21 (G.L.A.)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -2.857142857142
$ws.Range("G21").Value = 149
$ws.Range("H21").Value = 26.845637583892
$ws.Range("I21").Value = 1931
$ws.Range("J21").Value = 1365
$ws.Range("K21").Value = 41.465201465201
$ws.Range("L21").Value = 50.155520995334
$ws.Range("M21").Value = 27.206851119894
$ws.Range("N21").Value = -70.437844458052

# ---------------------------------------------------------------------
# Row 22 (TOTAL) — C22 flips from a numeric 4 to the literal text "0"
# ---------------------------------------------------------------------
$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 26
$ws.Range("K22").Value = 123.076923076923

# ---------------------------------------------------------------------
# Row 24 (Transit)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 46
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 187.5
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = 65.486725663716
$ws.Range("I24").Value = 1676
$ws.Range("J24").Value = 1185
$ws.Range("K24").Value = 41.434599156118
$ws.Range("L24").Value = 36.482084690553
$ws.Range("M24").Value = 68.273092369477

# ---------------------------------------------------------------------
# Row 25 (Housing)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 22.222222222222
$ws.Range("F25").Value = 81
$ws.Range("G25").Value = 88
$ws.Range("H25").Value = -7.954545454545
$ws.Range("I25").Value = 794
$ws.Range("J25").Value = 743
$ws.Range("K25").Value = 6.864064602960
$ws.Range("L25").Value = 19.758672699849
$ws.Range("M25").Value = -3.170731707317

# ---------------------------------------------------------------------
# Row 26 (Petit Larceny) — C26 flips from the literal text "0" to numeric 1
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 1
$ws.Range("D26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -42.857142857142
$ws.Range("I26").Value = 47
$ws.Range("J26").Value = 48
$ws.Range("K26").Value = -2.083333333333
$ws.Range("L26").Value = 17.5

# ---------------------------------------------------------------------
# Row 27 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 19
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 58.333333333333
$ws.Range("I27").Value = 102
$ws.Range("J27").Value = 94
$ws.Range("K27").Value = 8.510638297872
$ws.Range("L27").Value = 39.726027397260
